$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Great Deals" -> "Great deals" for boilers row
$ws.Range("B2").Value = "Great deals"

# Update percentages for showers / Bathroom / enclosures & trays / drainage & ufh
$ws.Range("C9").Value = "45%"
$ws.Range("C10").Value = "40%"
$ws.Range("C11").Value = "40%"
$ws.Range("C13").Value = "25%"

# fittings row becomes a "Great Deals" row (no percentage)
$ws.Range("B14").Value = "Great Deals"
$ws.Range("C14").Value = ""

# power tools / hand tools switch to "Save over" template
$ws.Range("B15").Value = "Save over"
$ws.Range("C15").Value = "30%"
$ws.Range("B16").Value = "Save over"
$ws.Range("C16").Value = "35%"

# New "test equipment" category row (was trade essentials)
$ws.Range("A17").Value = "test equipment"
$ws.Range("B17").Value = "Save over"
$ws.Range("C17").Value = "20%"

# trade essentials shifts down (was ventilation)
$ws.Range("A18").Value = "trade essentials"
$ws.Range("B18").Value = "Save over"
$ws.Range("C18").Value = "40%"

# ventilation shifts down (was kitchen)
$ws.Range("A19").Value = "ventilation"
$ws.Range("B19").Value = "Great Deals"
$ws.Range("C19").Value = ""

# New row: kitchen
$ws.Range("A20").Value = "kitchen"
$ws.Range("B20").Value = "Save"
$ws.Range("C20").Value = "20%"
